# Brazil and Nashville updates - Nashville_validation.xlsx
#
# Updates the "Validation" sheet's row-3 data (columns A-D and I-L),
# refreshes the current selection/view, and turns on portrait page
# orientation for the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Validation")

# --- Update cell values on row 3 (updated Nashville validation numbers) ---
$ws.Range("A3").Value = 2846
$ws.Range("B3").Value = 3716
$ws.Range("C3").Value = 4611
$ws.Range("D3").Value = 9083
$ws.Range("I3").Value = 34
$ws.Range("J3").Value = 44
$ws.Range("K3").Value = 52
$ws.Range("L3").Value = 111

# --- Update sheet selection to the refreshed data row ---
$ws.Activate()
$ws.Range("A3:P3").Select()

# --- Switch page setup to portrait orientation ---
$ws.PageSetup.Orientation = 1
